# TrapCardData.xlsx rule-change edit
# - Row 9  (传送阵 / Teleporter): new effect text
# - Row 13 (遥控炸弹 / Remote bomb -> 延迟爆弹 / Delayed bomb): renamed card
# - Row 14 (暗门 / Secret door): new effect text
# - View state: scroll sheet back to top-left, move selection to D15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - card renamed from "遥控炸弹"/"Remote bomb" to "延迟爆弹"/"Delayed bomb"
$ws.Range("A13").Value = "延迟爆弹"
$ws.Range("E13").Value = "Delayed bomb"

# Row 9 - Teleporter (传送阵) effect text changed
$ws.Range("D9").Value = "交锋时：如果同一行中怪物牌数量大于1，则将同一行中所有怪物牌洗回主牌堆，然后将本牌送墓。"

# Row 14 - Secret door (暗门) effect text changed
$ws.Range("D14").Value = "交锋时：选对位的1张怪物牌移动到本牌所在槽位，然后将本牌送墓。"

# View state: scroll back to top (removes topLeftCell="A11") and reselect D15
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D15").Select()
